# Added mouse control of modmatrix for testing.
# Rename the Morph X/Y modulation-matrix rows to Mouse X/Y, and repoint the
# FX parameter-path header row at the granular-synth parameter names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Row 3/4 labels: MorphX/MorphY -> Mousex/Mousey
$ws.Cells.Item(3, 1).Value = "Mousex"
$ws.Cells.Item(4, 1).Value = "Mousey"

# Row 2 header: OSC fx-param paths -> granular-synth parameter names
$ws.Cells.Item(2, 2).Value = "GrainRate"
$ws.Cells.Item(2, 3).Value = "Duration"
$ws.Cells.Item(2, 4).Value = "Freq"
$ws.Cells.Item(2, 5).Value = "FmPitch"
$ws.Cells.Item(2, 6).Value = "FmIndex"
$ws.Cells.Item(2, 7).Value = "Env"
$ws.Cells.Item(2, 8).Value = "Distr"
$ws.Cells.Item(2, 9).Value = "RndMask"

# Move the active selection to H3, matching the saved view state
$ws.Range("H3").Select()
